$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Move the existing "Kapasitans (nF)" column (D) one slot to the right (to E),
# carrying its values and header formatting along, so column D is freed up
# for the new "Admittans (p.u.)" data.
$ws.Range("D1:D10").Copy($ws.Range("E1"))
$ws.Range("D1:D10").ClearContents()

# New header for the freed-up column.
$ws.Range("D1").Value = "Admittans (p.u.)"

# New admittans (p.u.) values for each line.
$ws.Range("D2").Value = "(0.8429-9.0614j)"
$ws.Range("D3").Value = "(1.0268-11.0376j)"
$ws.Range("D4").Value = "(0.4038-4.3411j)"
$ws.Range("D5").Value = "(2.0322-22.3547j)"
$ws.Range("D6").Value = "(0.2459-3.9348j)"
$ws.Range("D7").Value = "(0.3219-5.1499j)"
$ws.Range("D8").Value = "(0.3907-4.2005j)"
$ws.Range("D9").Value = "(0.1824-1.9613j)"
$ws.Range("D10").Value = "(0.1654-1.7781j)"

# Widen the data columns (no longer auto "best fit") to comfortably fit the
# longer complex-number strings.
$ws.Columns.Item(2).ColumnWidth = 18
$ws.Columns.Item(3).ColumnWidth = 17.833333333333332
$ws.Columns.Item(4).ColumnWidth = 16.333333333333332

# Leave the selection where the author left it after finishing the edit.
$ws.Range("G9").Select()
